$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet: insert a new blank column before column N.
# This pushes the existing "Late" / "heading" / "Over Due" columns one
# column to the right (N->O, O->P, P->Q) and leaves the new column N blank.
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N").Insert() | Out-Null

# "Repayment Schedule" becomes the active sheet/tab (previously
# "Transactions" was active), with cell M19 selected.
$ws.Activate() | Out-Null
$ws.Range("M19").Select() | Out-Null
